# -----------------------------------------------------------------------
# "Actualizacion completa" update for BACKUP/SUPORT IMAGEN.xlsx
#
# Hoja2 holds the catalogue rows (cols A:L) plus a helper CONCATENATE
# column (O). For each product row the stock (H) is refreshed and the old
# numeric "source" flags in K/L are replaced by a single text label in K
# (K became a promo/status tag: "Oferta Flash", "Preventa 26/..",
# "Ultimas Unidades") while L is cleared out entirely.
#
# Afterwards Hoja2!O1:O24 (the rebuilt CSV strings) is copied as values
# into Hoja1!A1:A24, mirroring the manual copy/paste the author did, and
# the active sheet ends up back on Hoja1.
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# New stock (H) + status label (K) per data row on Hoja2.
# L is always cleared (previously a constant "1" helper flag).
$updates = @(
    @{ Row = 2;  Stock = 10; Label = "Oferta Flash" },
    @{ Row = 3;  Stock = 5;  Label = "Oferta Flash" },
    @{ Row = 4;  Stock = 10; Label = "Oferta Flash" },
    @{ Row = 5;  Stock = 5;  Label = "Oferta Flash" },
    @{ Row = 6;  Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 7;  Stock = 5;  Label = "Últimas Unidades" },
    @{ Row = 8;  Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 9;  Stock = 5;  Label = "Últimas Unidades" },
    @{ Row = 10; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 11; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 12; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 13; Stock = 10; Label = "Preventa 26/09" },
    @{ Row = 14; Stock = 10; Label = "Preventa 26/10" },
    @{ Row = 15; Stock = 10; Label = "Preventa 26/11" },
    @{ Row = 16; Stock = 10; Label = "Preventa 26/12" },
    @{ Row = 17; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 18; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 19; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 20; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 21; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 22; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 23; Stock = 2;  Label = "Últimas Unidades" },
    @{ Row = 24; Stock = 2;  Label = "Últimas Unidades" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws2.Cells.Item($r, 8).Value = $u.Stock        # H: stock
    $ws2.Cells.Item($r, 11).Value = $u.Label        # K: status / source label
    $ws2.Cells.Item($r, 12).ClearContents()         # L: no longer used
}

$wb.Application.Calculate()

# Copy the refreshed CSV helper column (Hoja2!O) into Hoja1!A as values,
# same as the manual copy -> paste values the author performed.
$ws2.Activate()
$ws2.Range("O1:O24").Select()
$ws2.Range("O1:O24").Copy()

$ws1.Activate()
$ws1.Range("A1:A24").Select()
$ws1.Range("A1").PasteSpecial(-4163)  # xlPasteValues

$wb.Application.CutCopyMode = $false
